$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25, columns B,C,D,E,F,I,J,K,L,M,N
# (column H / G / A are unchanged by this edit)
$newValues = @(
    @(1.02, 1.02037819947043, 1.023022984938293, 1.021392450590189, 1.031468292678807, 1.029715671014253, 1.025575648424696, 1.02585564403308, 1.024229915697476, 1.034276323527166, 1.01259007378955),
    @(1.02, 1.021267779614528, 1.023824780390868, 1.022145536458064, 1.03257997390409, 1.029830438728, 1.026102702764771, 1.026464243744226, 1.024789590829566, 1.035195762422167, 1.012766349669682),
    @(1.02, 1.021843813588545, 1.024344291939946, 1.022633584061396, 1.033299832761817, 1.02990331928466, 1.026443547113487, 1.026858094548972, 1.02515182775153, 1.035790645139107, 1.012880291868348),
    @(1.02, 1.022086076814971, 1.024562859891592, 1.022838937570022, 1.033602587035742, 1.029933626966795, 1.02658679043125, 1.027023679372222, 1.025304132570289, 1.036040719908201, 1.012928164051019),
    @(1.02, 1.022126759624015, 1.02459956806495, 1.022873427741769, 1.033653428161946, 1.029938696317892, 1.026610838792959, 1.027051482344123, 1.02532970638611, 1.036082707728859, 1.012936200282912),
    @(1.02, 1.021847050335772, 1.024347211807863, 1.022636327306234, 1.033303877686303, 1.029903725559802, 1.026445461326122, 1.026860307063777, 1.025153862777075, 1.035793986705858, 1.012880931653674),
    @(1.02, 1.020678750431262, 1.023293810459198, 1.021646803196521, 1.031843881181392, 1.029754743121663, 1.025753808722151, 1.026061312955064, 1.024419041324744, 1.034587063394048, 1.012649671690531),
    @(1.02, 1.018623291498596, 1.021442976816771, 1.019908947326655, 1.02927523705095, 1.029481655546126, 1.024533576495511, 1.024653784370856, 1.023124925432931, 1.032459909061177, 1.012241258678141),
    @(1.02, 1.017255223530733, 1.020212796640013, 1.018754368387555, 1.027565561875635, 1.029292525111905, 1.023719167971025, 1.023715765844001, 1.022262736734586, 1.031041568873264, 1.011968397678806),
    @(1.02, 1.016663378751254, 1.01968101187794, 1.018255387271605, 1.026825911843994, 1.029208958121526, 1.023366313355069, 1.023309685975716, 1.021889544130142, 1.03042736182162, 1.011850111460998),
    @(1.02, 1.016443623024017, 1.019483618618724, 1.018070188976194, 1.02655127093686, 1.029177666832016, 1.023235216739752, 1.023158864151213, 1.021750945952768, 1.030199209917015, 1.011806154674258),
    @(1.02, 1.016490757656948, 1.019525953999397, 1.018109908032195, 1.026610177901975, 1.029184390261307, 1.02326333878351, 1.02319121527977, 1.021780674712588, 1.03024814963299, 1.011815584454437),
    @(1.02, 1.016645212000317, 1.019664692522038, 1.018240075746311, 1.026803207939031, 1.029206376687175, 1.023355477489004, 1.023297218698493, 1.02187808710336, 1.030408502875986, 1.011846478385623),
    @(1.02, 1.016740387284824, 1.019750191871263, 1.018320295694195, 1.026922153117028, 1.029219890031374, 1.02341224314536, 1.023362532832067, 1.021938109060873, 1.030507300746005, 1.011865510514166),
    @(1.02, 1.017294513732736, 1.020248108339631, 1.018787504449456, 1.027614663745351, 1.029298035983574, 1.023742581415, 1.02374271797812, 1.022287507352063, 1.031082330619603, 1.011976245118737),
    @(1.02, 1.017642247267453, 1.020560677792003, 1.019080829949031, 1.028049231763005, 1.029346607479901, 1.023949738204661, 1.023981222467091, 1.022506714083413, 1.031443016861081, 1.012045669955724),
    @(1.02, 1.017845126166322, 1.020743080178659, 1.019252014203698, 1.028302770688478, 1.029374777032457, 1.024070548898951, 1.024120346554072, 1.022634587177117, 1.031653393574545, 1.012086151180967),
    @(1.02, 1.017914311385814, 1.020805289209792, 1.019310399249774, 1.028389231544927, 1.029384354735458, 1.024111738743428, 1.024167785674495, 1.022678190859033, 1.031725125657066, 1.012099951998584),
    @(1.02, 1.017604933383788, 1.020527133151624, 1.019049349325899, 1.028002600234295, 1.029341412906997, 1.023927514340648, 1.023955632310617, 1.022483193870283, 1.031404319203467, 1.01203822268096),
    @(1.02, 1.016599726780447, 1.019623833739801, 1.018201740570865, 1.026746362695224, 1.029199909152124, 1.023328345776958, 1.02326600293006, 1.021849400961967, 1.030361283061721, 1.011837381442673),
    @(1.02, 1.015968187116833, 1.019056676903777, 1.017669657862973, 1.025957083720356, 1.029109489363877, 1.022951447545468, 1.022832488931914, 1.021451038991794, 1.029705438427794, 1.011710988937483),
    @(1.02, 1.01630293301425, 1.019357262766423, 1.017951644641509, 1.026375441551112, 1.029157559969292, 1.023151264921867, 1.023062294653332, 1.021662205661634, 1.030053118337501, 1.011778002839771),
    @(1.02, 1.017621793761639, 1.020542290265381, 1.019063573778911, 1.028023670820896, 1.029343760609568, 1.023937556411157, 1.023967195374789, 1.022493821600518, 1.031421805026055, 1.012041587824396),
    @(1.02, 1.019154287071616, 1.021920813970251, 1.020357527310487, 1.029938809669431, 1.029553503844934, 1.024849202084467, 1.025017610098457, 1.023459392093332, 1.033009873658841, 1.012346948047582)
)

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowVals = $newValues[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $cellRef = "" + $cols[$j] + $rowIndex
        $ws.Range($cellRef).Value = $rowVals[$j]
    }
}

Write-Host "Updated $($newValues.Count) rows x $($cols.Count) columns of vm_pu data"